$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6:A10").Value = 4.707040960416666
$ws.Range("A11").Value = 4.256446442824075
$ws.Range("A12").Value = 7.63804360763889
$ws.Range("A37:A54").Value = 4.889766951157407
$ws.Range("A57:A61").Value = 5.769767546990741
$ws.Range("A62:A83").Value = 4.980157999768519
$ws.Range("A88:A89").Value = 6.685730172916668
$ws.Range("A105:A143").Value = 6.405517922222222
$ws.Range("A144:A148").Value = 5.827695573842593
$ws.Range("A149").Value = 8.284971370601852
